# permak dikit excel main_owner.py
# - "Metode bayar" (I) replaced by two numeric columns: "Pembulatan" (I) and "Bayar" (J)
# - Report sub-title reworded
# - TOTAL row now totals the new "Bayar" column (J) instead of the old text column (I)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Rows that previously had a non-zero "Pembulatan" (rounding) amount baked into
# Grand total via the payment-method note; everything else rounds to 0.
$roundingAmounts = @{ 41 = 500; 78 = 500 }

# --- Title text tweak -----------------------------------------------------
$ws.Range("A2").Value = "LAPORAN PENJUALAN PERIODE 01-06-2025 s/d 30-06-2025"

# --- Header row -------------------------------------------------------------
$ws.Range("I4").Value = "Pembulatan"

# J4 is a brand-new header cell; give it the same look as the other header
# cells (bold, filled, bordered) before writing its text.
$ws.Range("I4").Copy() | Out-Null
$ws.Range("J4").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("J4").Value = "Bayar"
$ws.Application.CutCopyMode = $false

# --- Data rows 5-85 -----------------------------------------------------
$firstRow = 5
$lastRow = 85

for ($r = $firstRow; $r -le $lastRow; $r++) {

    $grandTotal = $ws.Cells.Item($r, 8).Value2
    if ($grandTotal -eq $null) {
        continue
    }

    $rounding = 0
    if ($roundingAmounts.ContainsKey($r)) {
        $rounding = $roundingAmounts[$r]
    }
    $bayar = $grandTotal + $rounding

    # Give I and J the same numeric look (style) as the Grand total column (H)
    $ws.Range("H" + $r).Copy() | Out-Null
    $ws.Range("I" + $r).PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range("J" + $r).PasteSpecial($xlPasteFormats) | Out-Null

    $ws.Range("I" + $r).Value = $rounding
    $ws.Range("J" + $r).Value = $bayar
}

$ws.Application.CutCopyMode = $false

# --- TOTAL row (87) -------------------------------------------------------
$total = 0
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $v = $ws.Cells.Item($r, 10).Value2
    if ($v -ne $null) {
        $total = $total + $v
    }
}

# J87 takes over the bold total style that I87 used to have
$ws.Range("I87").Copy() | Out-Null
$ws.Range("J87").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("J87").Value = $total

# I87 becomes a blank cell styled like the rest of the TOTAL row (e.g. H87)
$ws.Range("H87").Copy() | Out-Null
$ws.Range("I87").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("I87").ClearContents() | Out-Null

$ws.Application.CutCopyMode = $false

# --- Column widths (I narrower, new J column) ------------------------------
$ws.Columns.Item(9).ColumnWidth = 14.4
$ws.Columns.Item(10).ColumnWidth = 10.8
